$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 "Tipo", copying the formatting from the existing header cells (A1:C1)
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Tipo"

# Update existing numeric values in row 2
$ws.Range("B2").Value = 0.06147519137037851
$ws.Range("C2").Value = 0.9994168802598677

# Add new value cell D2 "single"
$ws.Range("D2").Value = "single"
